# Update column F ("dSF") values for the rows that were re-pulled / recalculated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 0
    12 = -1
    16 = -5
    21 = -2
    25 = 2
    27 = 2
    28 = 2
    34 = -1
    41 = 4
    42 = -1
    43 = -7
    45 = -3
    54 = -2
    56 = 0
    62 = -1
    64 = -2
    67 = 0
    71 = -2
    73 = -6
    75 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
